$d = $word.ActiveDocument

# Split the last bullet paragraph ("Created a user interface using Thymeleaf, ...")
# into two paragraphs: the text stays in the first paragraph, and a new
# (initially empty) paragraph - carrying the same list/indent formatting - is
# created right after it, taking over the paragraph's trailing empty run.
$target = "Created a user interface using Thymeleaf, incorporating inbox management, message composition, and reply functionalities"

$d.Content.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target + "^p", 2)
